# Fix cpu name in diagrams: "Xeon E5-2666G" -> "Xeon E5-2276G"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the CPU name. All chart series reference cell A4 (Tabelle1!$A$4),
# so updating this single cell propagates to every chart's cached strings.
$ws.Range("A4").Value = "Xeon E5-2276G"

# Update the view/selection state recorded on the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A5").Select()
